# add admin lock overlay & JWT auth for settings protection
#
# Extends the Users seed-data sheet with three new columns:
#   I: IsActive (boolean)
#   J: Gender   (text)
#   K: TimeZone (number)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ------------------------------------------------
# I1 already carries the "blank but styled" cell (style index 2) - just
# give it the new header text, keeping its existing formatting.
$ws.Range("I1").Value = "IsActive"
$ws.Range("J1").Value = "Gender"
$ws.Range("K1").Value = "TimeZone"

# ---- Row 2 (admin) -------------------------------------------------------
$ws.Range("I2").Value = $true
$ws.Range("J2").Value = "Male"
$ws.Range("K2").Value = 5

# ---- Row 3 (manager) ------------------------------------------------------
$ws.Range("I3").Value = $true
$ws.Range("J3").Value = "Male"
$ws.Range("K3").Value = 5
